$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-24 22:53:28"
$wsZhCn.Range("H3").Value = "2016-03-24 22:54:00"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-24 22:53:32"
$wsDeDe.Range("H3").Value = "2016-03-24 22:54:09"
